$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.261.95"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "3.500.26"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.19"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.25"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.26"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").Value = "4.094.43"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000182"
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.120"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "3.499.48"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.89"
$ws.Range("E16").Value = "  -4.45%  "
$ws.Range("D17").Value = "64.277.26"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.91"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.63"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.59"
$ws.Range("E21").Value = "  +4.46%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "3.639.94"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.36"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.28"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("D32").Value = "3.520.35"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E33").Value = "  +5.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.45"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.16"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.90"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.55"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.31"
$ws.Range("E39").Value = "  +3.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0782"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.17"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.40"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.18"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  +3.70%  "
$ws.Range("D47").Value = "2.460.28"
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.894"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  +0.78%  "
